# Organize esp32 <-> escon wiring
# Remapped AO1 of M3 to an ADC1 pin.
# Removed all three unused AO2 DI3 wirings and relevant diodes and capacitors.
# Update the BOM accordingly: fewer capacitors (row 5) and fewer diodes (row 9)
# are now needed, and their designator lists shrink to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: 0.1uF MLCC capacitors - quantity drops from 9 to 6, designators
# trimmed down to the six that remain in use.
$ws.Range("A5").Value = 6
$ws.Range("C5").Value = "C6, C7, C8, C9, C10, C11"

# Row 9: 3.3V bidirectional ESD/TVS diodes - quantity drops from 6 to 3,
# designators trimmed down to the three that remain in use.
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "D5, D6, D7"

# Restore the last active selection recorded in the sheet view.
$ws.Range("I22").Select()
